# Updates crypto price/volume/hour data (columns D, E, G) for rows 2-51,
# matching the "Updated symbol list" GitHub Actions commit.
#
# Values in D/E/G are stored as *text* (not numbers) in the source sheet
# (e.g. "333.10" keeps its trailing zero, "1.39%" keeps the percent sign
# as a literal character, "7" is a plain digit string) so we prefix each
# literal with a leading apostrophe to force Excel to store it as text
# instead of auto-converting to a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'333.10"
$ws.Range("E2").Value = "'1.39%"
$ws.Range("G2").Value = "'7"
$ws.Range("D3").Value = "'44.25"
$ws.Range("E3").Value = "'6.55%"
$ws.Range("G3").Value = "'7"
$ws.Range("D4").Value = "'5.805"
$ws.Range("E4").Value = "'3.38%"
$ws.Range("G4").Value = "'7"
$ws.Range("D5").Value = "'0.08344"
$ws.Range("E5").Value = "'2.14%"
$ws.Range("G5").Value = "'7"
$ws.Range("D6").Value = "'8.816"
$ws.Range("E6").Value = "'0.94%"
$ws.Range("G6").Value = "'7"
$ws.Range("D7").Value = "'4.498"
$ws.Range("E7").Value = "'-0.58%"
$ws.Range("G7").Value = "'7"
$ws.Range("D8").Value = "'1.981"
$ws.Range("E8").Value = "'-1.77%"
$ws.Range("G8").Value = "'7"
$ws.Range("E9").Value = "'-2.81%"
$ws.Range("G9").Value = "'7"
$ws.Range("D10").Value = "'0.9345"
$ws.Range("E10").Value = "'1.57%"
$ws.Range("G10").Value = "'7"
$ws.Range("D11").Value = "'0.1245"
$ws.Range("E11").Value = "'-2.37%"
$ws.Range("G11").Value = "'7"
$ws.Range("D12").Value = "'0.1951"
$ws.Range("E12").Value = "'0.15%"
$ws.Range("G12").Value = "'7"
$ws.Range("D13").Value = "'0.09606"
$ws.Range("E13").Value = "'2.05%"
$ws.Range("G13").Value = "'7"
$ws.Range("D14").Value = "'0.03938"
$ws.Range("E14").Value = "'4.46%"
$ws.Range("G14").Value = "'7"
$ws.Range("E15").Value = "'0.85%"
$ws.Range("G15").Value = "'7"
$ws.Range("D16").Value = "'0.001316"
$ws.Range("E16").Value = "'1.25%"
$ws.Range("G16").Value = "'7"
$ws.Range("D17").Value = "'0.006065"
$ws.Range("E17").Value = "'-3.33%"
$ws.Range("G17").Value = "'7"
$ws.Range("D18").Value = "'3.512"
$ws.Range("E18").Value = "'2.10%"
$ws.Range("G18").Value = "'7"
$ws.Range("E19").Value = "'0.38%"
$ws.Range("G19").Value = "'7"
$ws.Range("D20").Value = "'8.988"
$ws.Range("E20").Value = "'8.67%"
$ws.Range("G20").Value = "'7"
$ws.Range("E21").Value = "'-1.61%"
$ws.Range("G21").Value = "'7"
$ws.Range("D22").Value = "'0.2572"
$ws.Range("E22").Value = "'6.65%"
$ws.Range("G22").Value = "'7"
$ws.Range("D23").Value = "'0.04418"
$ws.Range("E23").Value = "'0.15%"
$ws.Range("G23").Value = "'7"
$ws.Range("D24").Value = "'0.001257"
$ws.Range("E24").Value = "'-0.14%"
$ws.Range("G24").Value = "'7"
$ws.Range("D25").Value = "'0.004403"
$ws.Range("E25").Value = "'0.69%"
$ws.Range("G25").Value = "'7"
$ws.Range("D26").Value = "'0.0001192"
$ws.Range("E26").Value = "'0.89%"
$ws.Range("G26").Value = "'7"
$ws.Range("G27").Value = "'7"
$ws.Range("G28").Value = "'7"
$ws.Range("G29").Value = "'7"
$ws.Range("G30").Value = "'7"
$ws.Range("G31").Value = "'7"
$ws.Range("G32").Value = "'7"
$ws.Range("G33").Value = "'7"
$ws.Range("G34").Value = "'7"
$ws.Range("G35").Value = "'7"
$ws.Range("G36").Value = "'7"
$ws.Range("G37").Value = "'7"
$ws.Range("G38").Value = "'7"
$ws.Range("D39").Value = "'0.02799"
$ws.Range("E39").Value = "'0.94%"
$ws.Range("G39").Value = "'7"
$ws.Range("D40").Value = "'0.05722"
$ws.Range("E40").Value = "'5.68%"
$ws.Range("G40").Value = "'7"
$ws.Range("D41").Value = "'0.007933"
$ws.Range("E41").Value = "'3.72%"
$ws.Range("G41").Value = "'7"
$ws.Range("E42").Value = "'0.60%"
$ws.Range("G42").Value = "'7"
$ws.Range("D43").Value = "'0.009030"
$ws.Range("E43").Value = "'0.78%"
$ws.Range("G43").Value = "'7"
$ws.Range("D44").Value = "'0.002113"
$ws.Range("E44").Value = "'-0.90%"
$ws.Range("G44").Value = "'7"
$ws.Range("D45").Value = "'0.01050"
$ws.Range("E45").Value = "'-10.18%"
$ws.Range("G45").Value = "'7"
$ws.Range("D46").Value = "'0.00007229"
$ws.Range("E46").Value = "'8.52%"
$ws.Range("G46").Value = "'7"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("G47").Value = "'7"
$ws.Range("E48").Value = "'0.95%"
$ws.Range("G48").Value = "'7"
$ws.Range("G49").Value = "'7"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("G50").Value = "'7"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.03%"
$ws.Range("G51").Value = "'7"
